$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Pepino ensalada at Vega
# Monumental Concepción. Insert a row above the current row 199 (shifting
# the existing rows 199:242 down to 200:243) and fill it with the new
# reading.
$ws.Rows("199:199").Insert()

$ws.Cells.Item(199, 1).Value = 11
$ws.Cells.Item(199, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(199, 3).Value = "Bíobío"
$ws.Cells.Item(199, 4).Value = 45204
$ws.Cells.Item(199, 5).Value = 8
$ws.Cells.Item(199, 6).Value = 100112043
$ws.Cells.Item(199, 7).Value = "Pepino ensalada"
$ws.Cells.Item(199, 8).Value = "Sin especificar"
$ws.Cells.Item(199, 9).Value = "Primera"
$ws.Cells.Item(199, 10).Value = 100
$ws.Cells.Item(199, 11).Value = 14000
$ws.Cells.Item(199, 12).Value = 15000
$ws.Cells.Item(199, 13).Value = 14500
$ws.Cells.Item(199, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(199, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(199, 16).Value = 242
$ws.Cells.Item(199, 17).Value = 60
$ws.Cells.Item(199, 18).Value = "Hortaliza"
